# DCIT 23 - IT 1B POINTS workbook update
# Adds Finals "Score" (AN), and computed "Excess"/"Items" columns (AO/AP/AQ)
# to the Recitation sheet, and updates scroll/selection state on a couple
# of sheets.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Recitation
$ws2 = $wb.Worksheets.Item(2)   # Quizzes&Ass

# ---------------------------------------------------------------------------
# Header additions on the Recitation sheet (row 1 / row 2)
# ---------------------------------------------------------------------------
$ws1.Range("AO1").Value = "Finals"
$ws1.Range("AP1").Value = 70

$ws1.Range("AO2").Value = "Score"
$ws1.Range("AP2").Value = "Excess"
$ws1.Range("AQ2").Value = "Items"

# ---------------------------------------------------------------------------
# Finals scores entered in column AN (rows 3-42; some rows stay blank)
# ---------------------------------------------------------------------------
$ws1.Range("AN3").Value = 33
$ws1.Range("AN4").Value = 48
$ws1.Range("AN5").Value = 68
$ws1.Range("AN6").Value = 39
$ws1.Range("AN8").Value = 33
$ws1.Range("AN9").Value = 41
$ws1.Range("AN10").Value = 20
$ws1.Range("AN11").Value = 25
$ws1.Range("AN12").Value = 36
$ws1.Range("AN15").Value = 53
$ws1.Range("AN16").Value = 50
$ws1.Range("AN17").Value = 51
$ws1.Range("AN18").Value = 45
$ws1.Range("AN19").Value = 10
$ws1.Range("AN20").Value = 14
$ws1.Range("AN21").Value = 46
$ws1.Range("AN22").Value = 46
$ws1.Range("AN23").Value = 37
$ws1.Range("AN24").Value = 59
$ws1.Range("AN25").Value = 47
$ws1.Range("AN26").Value = 53
$ws1.Range("AN27").Value = 44
$ws1.Range("AN31").Value = 60
$ws1.Range("AN32").Value = 24
$ws1.Range("AN33").Value = 64
$ws1.Range("AN34").Value = 51
$ws1.Range("AN35").Value = 47
$ws1.Range("AN36").Value = 28
$ws1.Range("AN37").Value = 61
$ws1.Range("AN38").Value = 37
$ws1.Range("AN39").Value = 65
$ws1.Range("AN40").Value = 58
$ws1.Range("AN41").Value = 52
$ws1.Range("AN42").Value = 46

# ---------------------------------------------------------------------------
# Items (max finals score), Score and Excess formulas for rows 3-59
# ---------------------------------------------------------------------------
$ws1.Range("AQ3:AQ59").Value = 70

$ws1.Range("AO3:AO59").Formula = "=IF(AN3+AL3>AQ3,AQ3,AN3+AL3)"
$ws1.Range("AP3:AP59").Formula = "=IF(AN3+AL3>AQ3,AL3+AN3-AQ3,0)"

$ws1.Range("AO3:AP59").NumberFormat = "0"

# ---------------------------------------------------------------------------
# View / selection state
# ---------------------------------------------------------------------------
# Update the Quizzes&Ass sheet's frozen-pane scroll position and selection
# first (without leaving it as the active/selected tab)...
$ws2.Range("R1").Select()
$ws2.Range("B11").Select()

# ...then return to the Recitation sheet so it stays the active tab, and
# move the selection/scroll position there too.
$ws1.Select()
$ws1.Range("AL1").Select()
$ws1.Range("AN7").Select()
